$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Sending cluster)
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "FAPs"
$ws.Range("A4").Value = "sCs"
$ws.Range("A5").Value = "sCs"

# Column B (Ligand symbol)
$ws.Range("B2").Value = "Ccl21b"
$ws.Range("B3").Value = "Ccl21b"
$ws.Range("B4").Value = "Ccl21b"
$ws.Range("B5").Value = "Ccl21b"

# Column C (Receptor symbol)
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("C5").Value = "Cxcr3"

# Column D (Target cluster)
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "M2"
$ws.Range("D4").Value = "ECs"
$ws.Range("D5").Value = "M2"

# Column E (Ligand-expressing cells)
$ws.Range("E2").Value = 3
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2

# Column F (Ligand detection rate)
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("F5").Value = 0.6666666666666666

# Column G (Ligand average expression value)
$ws.Range("G2").Value = 0.2073196666666667
$ws.Range("G3").Value = 0.2073196666666667
$ws.Range("G4").Value = 0.265412
$ws.Range("G5").Value = 0.265412

# Column H (Ligand total expression value)
$ws.Range("H2").Value = 0.6219589999999999
$ws.Range("H3").Value = 0.6219589999999999
$ws.Range("H4").Value = 0.7962360000000001
$ws.Range("H5").Value = 0.7962360000000001

# Column I (Ligand derived specificity of average expression value)
$ws.Range("I2").Value = 0.4385567570045022
$ws.Range("I3").Value = 0.4385567570045022
$ws.Range("I4").Value = 0.5614432429954979
$ws.Range("I5").Value = 0.5614432429954979

# Column J (Ligand derived specificity of total expression value)
$ws.Range("J2").Value = 0.4385567570045022
$ws.Range("J3").Value = 0.4385567570045022
$ws.Range("J4").Value = 0.5614432429954979
$ws.Range("J5").Value = 0.5614432429954979

# Column K (Receptor-expressing cells)
$ws.Range("K2").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("K4").Value = 1
$ws.Range("K5").Value = 3

# Column L (Receptor detection rate)
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("L3").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("L5").Value = 1

# Column M (Receptor average expression value)
$ws.Range("M2").Value = 1.863797
$ws.Range("M3").Value = 1.552879
$ws.Range("M4").Value = 1.863797
$ws.Range("M5").Value = 1.552879

# Column N (Receptor total expression value)
$ws.Range("N2").Value = 5.591391
$ws.Range("N3").Value = 4.658637000000001
$ws.Range("N4").Value = 5.591391
$ws.Range("N5").Value = 4.658637000000001

# Column O (Receptor derived specificity of average expression value)
$ws.Range("O2").Value = 0.5455000708290748
$ws.Range("O3").Value = 0.4544999291709252
$ws.Range("O4").Value = 0.5455000708290748
$ws.Range("O5").Value = 0.4544999291709252

# Column P (Receptor derived specificity of total expression value)
$ws.Range("P2").Value = 0.5455000708290748
$ws.Range("P3").Value = 0.4544999291709252
$ws.Range("P4").Value = 0.5455000708290748
$ws.Range("P5").Value = 0.4544999291709252

# Column Q (Edge average expression weight)
$ws.Range("Q2").Value = 0.3864017727743333
$ws.Range("Q3").Value = 0.3219423566536667
$ws.Range("Q4").Value = 0.494674089364
$ws.Range("Q5").Value = 0.4121527211480001

# Column R (Edge total expression weight)
$ws.Range("R2").Value = 3.477615954969
$ws.Range("R3").Value = 2.897481209883
$ws.Range("R4").Value = 4.452066804276
$ws.Range("R5").Value = 3.709374490332001

# Column S (Edge average expression derived specificity)
$ws.Range("S2").Value = 0.2392327420085253
$ws.Range("S3").Value = 0.1993240149959769
$ws.Range("S4").Value = 0.3062673288205495
$ws.Range("S5").Value = 0.2551759141749483

# Column T (Edge total expression derived specificity)
$ws.Range("T2").Value = 0.2392327420085253
$ws.Range("T3").Value = 0.1993240149959769
$ws.Range("T4").Value = 0.3062673288205495
$ws.Range("T5").Value = 0.2551759141749483
